$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "최종점수" (K) column values
$ws.Range("K2").Value = 58.5
$ws.Range("K3").Value = 56.7
$ws.Range("K4").Value = 55.5
$ws.Range("K5").Value = 55.5

# Update "MACRO_SCORE" (N) column values
$ws.Range("N2").Value = 54.83846622768671
$ws.Range("N3").Value = 54.83846622768671
$ws.Range("N4").Value = 54.83846622768671
$ws.Range("N5").Value = 54.83846622768671
